$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 883.3333  # H18: 975 -> 883.3333
$ws.Cells.Item(18, 10).Value = 700  # J18: 0 -> 700
$ws.Cells.Item(18, 12).Value = 700  # L18: 0 -> 700
$ws.Cells.Item(18, 14).Value = -1268  # N18: None -> -1268

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(39, 8).Value = 490  # H39: 674.1818 -> 490
$ws.Cells.Item(39, 9).Value = 490  # I39: 739.6 -> 490
$ws.Cells.Item(39, 10).Value = 0  # J39: 20 -> 0
$ws.Cells.Item(39, 11).Value = 1470  # K39: 2218.8 -> 1470
$ws.Cells.Item(39, 12).Value = 0  # L39: 60 -> 0
$ws.Cells.Item(39, 13).Value = -1174  # M39: -1922.8 -> -1174
$ws.Cells.Item(39, 14).ClearContents()  # N39: -652 -> (removed)

# ALC row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(99, 8).Value = 846.5714  # H99: 776.5 -> 846.5714
$ws.Cells.Item(99, 9).Value = 737.6667  # I99: 673.4286 -> 737.6667
$ws.Cells.Item(99, 10).Value = 1500  # J99: 1498 -> 1500
$ws.Cells.Item(99, 11).Value = 2213.0001  # K99: 2020.2858 -> 2213.0001
$ws.Cells.Item(99, 12).Value = 4500  # L99: 4494 -> 4500
$ws.Cells.Item(99, 13).Value = -715.0001000000002  # M99: -522.2857999999999 -> -715.0001000000002
$ws.Cells.Item(99, 14).Value = -7496  # N99: -7490 -> -7496

# ALC row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(101, 8).Value = 33334000  # H101: 50000500 -> 33334000
$ws.Cells.Item(101, 9).Value = 33334000  # I101: 50000500 -> 33334000
$ws.Cells.Item(101, 11).Value = 100002000  # K101: 150001500 -> 100002000
$ws.Cells.Item(101, 13).Value = -100000378  # M101: -149999878 -> -100000378

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 3577.5  # H112: 3781.75 -> 3577.5
$ws.Cells.Item(112, 10).Value = 3555  # J112: 3842.3333 -> 3555
$ws.Cells.Item(112, 12).Value = 10665  # L112: 11526.9999 -> 10665
$ws.Cells.Item(112, 14).Value = -12881  # N112: -13742.9999 -> -12881

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1649.75  # H132: 1579.8 -> 1649.75
$ws.Cells.Item(132, 9).Value = 1649.75  # I132: 1579.8 -> 1649.75
$ws.Cells.Item(132, 11).Value = 4949.25  # K132: 4739.4 -> 4949.25
$ws.Cells.Item(132, 13).Value = -2419.25  # M132: -2209.4 -> -2419.25

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6006  # H32: 5934.143 -> 6006
$ws.Cells.Item(32, 10).Value = 0  # J32: 5000 -> 0
$ws.Cells.Item(32, 12).Value = 0  # L32: 5000 -> 0
$ws.Cells.Item(32, 14).ClearContents()  # N32: -5574 -> (removed)

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 18643.428  # H63: 16363 -> 18643.428
$ws.Cells.Item(63, 9).Value = 18643.428  # I63: 16363 -> 18643.428
$ws.Cells.Item(63, 11).Value = 18643.428  # K63: 16363 -> 18643.428
$ws.Cells.Item(63, 13).Value = -17957.428  # M63: -15677 -> -17957.428

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 18643.428  # H66: 16363 -> 18643.428
$ws.Cells.Item(66, 9).Value = 18643.428  # I66: 16363 -> 18643.428
$ws.Cells.Item(66, 11).Value = 93217.14  # K66: 81815 -> 93217.14
$ws.Cells.Item(66, 13).Value = -89785.14  # M66: -78383 -> -89785.14

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2799.4285  # H132: 3000 -> 2799.4285
$ws.Cells.Item(132, 9).Value = 2719.4  # I132: 3000 -> 2719.4
$ws.Cells.Item(132, 10).Value = 2999.5  # J132: 3000 -> 2999.5
$ws.Cells.Item(132, 11).Value = 8158.200000000001  # K132: 9000 -> 8158.200000000001
$ws.Cells.Item(132, 12).Value = 8998.5  # L132: 9000 -> 8998.5
$ws.Cells.Item(132, 13).Value = -5628.200000000001  # M132: -6470 -> -5628.200000000001
$ws.Cells.Item(132, 14).Value = -14058.5  # N132: -14060 -> -14058.5

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 864.5  # H22: 795.3333 -> 864.5
$ws.Cells.Item(22, 9).Value = 527.4  # I22: 479.83334 -> 527.4
$ws.Cells.Item(22, 11).Value = 527.4  # K22: 479.83334 -> 527.4
$ws.Cells.Item(22, 13).Value = -354.4  # M22: -306.83334 -> -354.4

# BSM row 76
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(76, 8).Value = 141876  # H76: 112656.75 -> 141876
$ws.Cells.Item(76, 10).Value = 141876  # J76: 112656.75 -> 141876
$ws.Cells.Item(76, 12).Value = 141876  # L76: 112656.75 -> 141876
$ws.Cells.Item(76, 14).Value = -142506  # N76: -113286.75 -> -142506

# BSM row 79
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(79, 8).Value = 141876  # H79: 112656.75 -> 141876
$ws.Cells.Item(79, 10).Value = 141876  # J79: 112656.75 -> 141876
$ws.Cells.Item(79, 12).Value = 141876  # L79: 112656.75 -> 141876
$ws.Cells.Item(79, 14).Value = -144060  # N79: -114840.75 -> -144060

# BSM row 88
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(88, 8).Value = 31430  # H88: 42145 -> 31430
$ws.Cells.Item(88, 10).Value = 31430  # J88: 42145 -> 31430
$ws.Cells.Item(88, 12).Value = 31430  # L88: 42145 -> 31430
$ws.Cells.Item(88, 14).Value = -32242  # N88: -42957 -> -32242

# BSM row 91
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(91, 8).Value = 31430  # H91: 42145 -> 31430
$ws.Cells.Item(91, 10).Value = 31430  # J91: 42145 -> 31430
$ws.Cells.Item(91, 12).Value = 31430  # L91: 42145 -> 31430
$ws.Cells.Item(91, 14).Value = -34238  # N91: -44953 -> -34238

# CRP row 17
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value = 15000  # H17: 0 -> 15000
$ws.Cells.Item(17, 10).Value = 15000  # J17: 0 -> 15000
$ws.Cells.Item(17, 12).Value = 15000  # L17: 0 -> 15000
$ws.Cells.Item(17, 14).Value = -15348  # N17: None -> -15348

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2325.1667  # H31: 2290.2 -> 2325.1667
$ws.Cells.Item(31, 10).Value = 3147  # J31: 3470.5 -> 3147
$ws.Cells.Item(31, 12).Value = 3147  # L31: 3470.5 -> 3147
$ws.Cells.Item(31, 14).Value = -3737  # N31: -4060.5 -> -3737

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2325.1667  # H34: 2290.2 -> 2325.1667
$ws.Cells.Item(34, 10).Value = 3147  # J34: 3470.5 -> 3147
$ws.Cells.Item(34, 12).Value = 3147  # L34: 3470.5 -> 3147
$ws.Cells.Item(34, 14).Value = -3551  # N34: -3874.5 -> -3551

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2939.6  # H58: 2981.5 -> 2939.6
$ws.Cells.Item(58, 9).Value = 2939.6  # I58: 2981.5 -> 2939.6
$ws.Cells.Item(58, 11).Value = 2939.6  # K58: 2981.5 -> 2939.6
$ws.Cells.Item(58, 13).Value = -2736.6  # M58: -2778.5 -> -2736.6

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 2805.9  # H105: 2868.3 -> 2805.9
$ws.Cells.Item(105, 9).Value = 2235.2856  # I105: 2369 -> 2235.2856
$ws.Cells.Item(105, 11).Value = 2235.2856  # K105: 2369 -> 2235.2856
$ws.Cells.Item(105, 13).Value = -488.2856000000002  # M105: -622 -> -488.2856000000002

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 2954.182  # H122: 3298.7 -> 2954.182
$ws.Cells.Item(122, 9).Value = 2954.182  # I122: 3298.7 -> 2954.182
$ws.Cells.Item(122, 11).Value = 8862.545999999998  # K122: 9896.099999999999 -> 8862.545999999998
$ws.Cells.Item(122, 13).Value = -6412.545999999998  # M122: -7446.099999999999 -> -6412.545999999998

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1381.2222  # H134: 932.1111 -> 1381.2222
$ws.Cells.Item(134, 9).Value = 928.875  # I134: 932.1111 -> 928.875
$ws.Cells.Item(134, 10).Value = 5000  # J134: 0 -> 5000
$ws.Cells.Item(134, 11).Value = 2786.625  # K134: 2796.3333 -> 2786.625
$ws.Cells.Item(134, 12).Value = 15000  # L134: 0 -> 15000
$ws.Cells.Item(134, 13).Value = -251.625  # M134: -261.3332999999998 -> -251.625
$ws.Cells.Item(134, 14).Value = -20070  # N134: None -> -20070

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2939.6  # H136: 2981.5 -> 2939.6
$ws.Cells.Item(136, 9).Value = 2939.6  # I136: 2981.5 -> 2939.6
$ws.Cells.Item(136, 11).Value = 8818.799999999999  # K136: 8944.5 -> 8818.799999999999
$ws.Cells.Item(136, 13).Value = -6268.799999999999  # M136: -6394.5 -> -6268.799999999999

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 11001923  # H4: 12224461 -> 11001923
$ws.Cells.Item(4, 10).Value = 3178.5  # J4: 3999 -> 3178.5
$ws.Cells.Item(4, 12).Value = 9535.5  # L4: 11997 -> 9535.5
$ws.Cells.Item(4, 14).Value = -9759.5  # N4: -12221 -> -9759.5

# CUL row 51
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 1874.5  # H51: 1999.3334 -> 1874.5
$ws.Cells.Item(51, 9).Value = 1500  # I51: 0 -> 1500
$ws.Cells.Item(51, 11).Value = 4500  # K51: 0 -> 4500
$ws.Cells.Item(51, 13).Value = -4040  # M51: None -> -4040

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 2040.2858  # H113: 1706.5555 -> 2040.2858
$ws.Cells.Item(113, 9).Value = 2645  # I113: 1240.6666 -> 2645
$ws.Cells.Item(113, 11).Value = 7935  # K113: 3721.9998 -> 7935
$ws.Cells.Item(113, 13).Value = -5765  # M113: -1551.9998 -> -5765

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 17506.7  # H121: 17250.4 -> 17506.7
$ws.Cells.Item(121, 9).Value = 60015  # I121: 40486 -> 60015
$ws.Cells.Item(121, 10).Value = 6879.625  # J121: 7292.2856 -> 6879.625
$ws.Cells.Item(121, 11).Value = 180045  # K121: 121458 -> 180045
$ws.Cells.Item(121, 12).Value = 20638.875  # L121: 21876.8568 -> 20638.875
$ws.Cells.Item(121, 13).Value = -178735  # M121: -120148 -> -178735
$ws.Cells.Item(121, 14).Value = -23258.875  # N121: -24496.8568 -> -23258.875

# GSM row 19
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 0  # H19: 1055.2 -> 0
$ws.Cells.Item(19, 9).Value = 0  # I19: 10 -> 0
$ws.Cells.Item(19, 10).Value = 0  # J19: 1316.5 -> 0
$ws.Cells.Item(19, 11).Value = 0  # K19: 10 -> 0
$ws.Cells.Item(19, 12).Value = 0  # L19: 1316.5 -> 0
$ws.Cells.Item(19, 13).ClearContents()  # M19: 278 -> (removed)
$ws.Cells.Item(19, 14).ClearContents()  # N19: -1892.5 -> (removed)

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8632.5  # H70: 8343.416999999999 -> 8632.5
$ws.Cells.Item(70, 10).Value = 9972.362999999999  # J70: 9152.532999999999 -> 9972.362999999999
$ws.Cells.Item(70, 12).Value = 9972.362999999999  # L70: 9152.532999999999 -> 9972.362999999999
$ws.Cells.Item(70, 14).Value = -10512.363  # N70: -9692.532999999999 -> -10512.363

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 8632.5  # H73: 8343.416999999999 -> 8632.5
$ws.Cells.Item(73, 10).Value = 9972.362999999999  # J73: 9152.532999999999 -> 9972.362999999999
$ws.Cells.Item(73, 12).Value = 9972.362999999999  # L73: 9152.532999999999 -> 9972.362999999999
$ws.Cells.Item(73, 14).Value = -11844.363  # N73: -11024.533 -> -11844.363

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6069.2856  # H7: 6032.3794 -> 6069.2856
$ws.Cells.Item(7, 9).Value = 3829.1667  # I7: 3765.3076 -> 3829.1667
$ws.Cells.Item(7, 10).Value = 7749.375  # J7: 7874.375 -> 7749.375
$ws.Cells.Item(7, 11).Value = 3829.1667  # K7: 3765.3076 -> 3829.1667
$ws.Cells.Item(7, 12).Value = 7749.375  # L7: 7874.375 -> 7749.375
$ws.Cells.Item(7, 13).Value = -3717.1667  # M7: -3653.3076 -> -3717.1667
$ws.Cells.Item(7, 14).Value = -7973.375  # N7: -8098.375 -> -7973.375

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1718.6  # H16: 1682.1666 -> 1718.6
$ws.Cells.Item(16, 9).Value = 1773.5  # I16: 1718.8 -> 1773.5
$ws.Cells.Item(16, 11).Value = 1773.5  # K16: 1718.8 -> 1773.5
$ws.Cells.Item(16, 13).Value = -1603.5  # M16: -1548.8 -> -1603.5

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2732.75  # H46: 2979.5 -> 2732.75
$ws.Cells.Item(46, 9).Value = 2578.4  # I46: 2973.25 -> 2578.4
$ws.Cells.Item(46, 10).Value = 2843  # J46: 2983.6667 -> 2843
$ws.Cells.Item(46, 11).Value = 2578.4  # K46: 2973.25 -> 2578.4
$ws.Cells.Item(46, 12).Value = 2843  # L46: 2983.6667 -> 2843
$ws.Cells.Item(46, 13).Value = -2390.4  # M46: -2785.25 -> -2390.4
$ws.Cells.Item(46, 14).Value = -3219  # N46: -3359.6667 -> -3219

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 181.23077  # H55: 195.08333 -> 181.23077
$ws.Cells.Item(55, 9).Value = 94.875  # I55: 106.28571 -> 94.875
$ws.Cells.Item(55, 11).Value = 94.875  # K55: 106.28571 -> 94.875
$ws.Cells.Item(55, 13).Value = 78.125  # M55: 66.71429000000001 -> 78.125

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 6069.2856  # H126: 6032.3794 -> 6069.2856
$ws.Cells.Item(126, 9).Value = 3829.1667  # I126: 3765.3076 -> 3829.1667
$ws.Cells.Item(126, 10).Value = 7749.375  # J126: 7874.375 -> 7749.375
$ws.Cells.Item(126, 11).Value = 11487.5001  # K126: 11295.9228 -> 11487.5001
$ws.Cells.Item(126, 12).Value = 23248.125  # L126: 23623.125 -> 23248.125
$ws.Cells.Item(126, 13).Value = -9017.500100000001  # M126: -8825.9228 -> -9017.500100000001
$ws.Cells.Item(126, 14).Value = -28188.125  # N126: -28563.125 -> -28188.125

# WVR row 15
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 39999  # H15: 0 -> 39999
$ws.Cells.Item(15, 10).Value = 39999  # J15: 0 -> 39999
$ws.Cells.Item(15, 12).Value = 39999  # L15: 0 -> 39999
$ws.Cells.Item(15, 14).Value = -40575  # N15: None -> -40575

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 739  # H122: 748.75 -> 739
$ws.Cells.Item(122, 9).Value = 739  # I122: 748.75 -> 739
$ws.Cells.Item(122, 11).Value = 2217  # K122: 2246.25 -> 2217
$ws.Cells.Item(122, 13).Value = 233  # M122: 203.75 -> 233

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4935.0713  # H132: 5890.5454 -> 4935.0713
$ws.Cells.Item(132, 9).Value = 5062.636  # I132: 6424.25 -> 5062.636
$ws.Cells.Item(132, 11).Value = 15187.908  # K132: 19272.75 -> 15187.908
$ws.Cells.Item(132, 13).Value = -12657.908  # M132: -16742.75 -> -12657.908
